$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.017.92'
$ws.Range('D3').Value = '1.640.86'
$ws.Range('E3').Value = '  -0.50%  '
$c = $ws.Range('D4')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '1.002'
$c.Style = $origStyle
$ws.Range('E4').Value = '  -0.71%  '
$c = $ws.Range('D5')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '215.84'
$c.Style = $origStyle
$ws.Range('E5').Value = '  +0.13%  '
$c = $ws.Range('D6')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.5160'
$c.Style = $origStyle
$ws.Range('E6').Value = '  +1.57%  '
$c = $ws.Range('D7')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.Style = $origStyle
$ws.Range('E7').Value = '  -0.60%  '
$ws.Range('E8').Value = '  +0.51%  '
$c = $ws.Range('D9')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.06383'
$c.Style = $origStyle
$ws.Range('E9').Value = '  -0.70%  '
$c = $ws.Range('D10')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '19.89'
$c.Style = $origStyle
$ws.Range('E10').Value = '  +0.89%  '
$c = $ws.Range('D11')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.07766'
$c.Style = $origStyle
$ws.Range('E11').Value = '  -0.23%  '
$c = $ws.Range('D12')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '4.297'
$c.Style = $origStyle
$ws.Range('E12').Value = '  -0.40%  '
$ws.Range('D13').Value = '1.632.74'
$ws.Range('E13').Value = '  -1.36%  '
$c = $ws.Range('D14')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.5482'
$c.Style = $origStyle
$ws.Range('E14').Value = '  +0.04%  '
$ws.Range('D15').Value = '0.0₅7798'
$ws.Range('E15').Value = '  -1.46%  '
$ws.Range('E16').Value = '  -0.93%  '
$ws.Range('D17').Value = '26.031.63'
$ws.Range('E17').Value = '  +0.15%  '
$ws.Range('E18').Value = '  -0.54%  '
$c = $ws.Range('D19')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '199.47'
$c.Style = $origStyle
$c = $ws.Range('D20')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '4.469'
$c.Style = $origStyle
$ws.Range('E20').Value = '  +0.82%  '
$c = $ws.Range('D21')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '9.998'
$c.Style = $origStyle
$ws.Range('E21').Value = '  -0.39%  '
$ws.Range('E22').Value = '  +0.78%  '
$c = $ws.Range('D23')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '1.002'
$c.Style = $origStyle
$ws.Range('E23').Value = '  -0.86%  '
$c = $ws.Range('D24')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '1.897'
$c.Style = $origStyle
$ws.Range('E24').Value = '  +2.25%  '
$c = $ws.Range('D25')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '142.31'
$c.Style = $origStyle
$ws.Range('E25').Value = '  +0.64%  '
$c = $ws.Range('D26')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.1233'
$c.Style = $origStyle
$ws.Range('E26').Value = '  +7.44%  '
$c = $ws.Range('D27')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '6.885'
$c.Style = $origStyle
$ws.Range('E27').Value = '  -0.37%  '
$ws.Range('E28').Value = '  -0.88%  '
$ws.Range('E29').Value = '  +0.02%  '
$c = $ws.Range('D30')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.04868'
$c.Style = $origStyle
$ws.Range('E30').Value = '  -3.41%  '
$c = $ws.Range('D31')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.314'
$c.Style = $origStyle
$ws.Range('E31').Value = '  +1.03%  '
$c = $ws.Range('D32')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.245'
$c.Style = $origStyle
$ws.Range('E32').Value = '  +1.11%  '
$c = $ws.Range('D33')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '1.546'
$c.Style = $origStyle
$ws.Range('E33').Value = '  +0.03%  '
$c = $ws.Range('D34')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.378'
$c.Style = $origStyle
$ws.Range('E34').Value = '  +0.33%  '
$c = $ws.Range('D35')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.9215'
$c.Style = $origStyle
$ws.Range('E35').Value = '  +3.02%  '
$c = $ws.Range('D36')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.5600'
$c.Style = $origStyle
$ws.Range('E36').Value = '  +0.85%  '
$c = $ws.Range('D37')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.568'
$c.Style = $origStyle
$ws.Range('E37').Value = '  -1.11%  '
$ws.Range('D38').Value = '1.118.19'
$ws.Range('E38').Value = '  -1.55%  '
$c = $ws.Range('D39')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.01577'
$c.Style = $origStyle
$ws.Range('E39').Value = '  +0.78%  '
$ws.Range('E40').Value = '  -0.74%  '
$c = $ws.Range('D41')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.528'
$c.Style = $origStyle
$ws.Range('E41').Value = '  -1.22%  '
$c = $ws.Range('D42')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '5.580'
$c.Style = $origStyle
$ws.Range('E42').Value = '  -1.62%  '
$c = $ws.Range('D43')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.8089'
$c.Style = $origStyle
$ws.Range('E43').Value = '  -0.85%  '
$c = $ws.Range('D44')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '99.76'
$c.Style = $origStyle
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('D45').Value = '0.0₈121'
$ws.Range('E45').Value = '  -0.90%  '
$ws.Range('D46').Value = '1.785.25'
$ws.Range('E46').Value = '  -0.01%  '
$c = $ws.Range('D47')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.4535'
$c.Style = $origStyle
$ws.Range('E47').Value = '  -0.19%  '
$ws.Range('E48').Value = '  +0.26%  '
$c = $ws.Range('D49')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '1.008'
$c.Style = $origStyle
$ws.Range('E49').Value = '  +0.09%  '
$c = $ws.Range('D50')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.05221'
$c.Style = $origStyle
$ws.Range('E50').Value = '  +2.45%  '
$c = $ws.Range('D51')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '7.525'
$c.Style = $origStyle
$ws.Range('E51').Value = '  +1.59%  '
